# RPG OBJECTS.xlsx - "Add files via upload" commit
# Had a slow day so I fixed a bunch of stuff.
#  - Moved loading bar a bit.                         -> selection moved to M18
#  - Slotting equipment works properly now.            -> I14/J14 swap, E/F/H values on rows 17-18
#  - Removed glitch from map N. exit.                   -> F7 value
#  - Added new items with terrible 80's dad jokes.      -> B17 "inanimate bones", B18 "Shoulder pad"
#  - Added new locations.                               -> O20 updated text, O26/O27 new trigger rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7: fixed north exit value ---
$ws.Range("F7").Value = 4

# --- Row 14: trigger column moved from I to J ---
$ws.Range("I14").ClearContents()
$ws.Range("J14").Value = 13

# --- O20: stat line "3)" now has an extra 3 in the 8th slot ---
$ws.Range("O20").Value = "3)  0, 4, 0, 0, 0,13, 0, 3, 0, 0"

# --- O26 / O27: two new location/trigger rows, matching the style used by O18:O25 ---
$ws.Range("O25").Copy()
$ws.Range("O26:O27").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("O27").Value = "10)  0 ,18, 0, 0, 0, 0, 0, 2, 0, 0"
$ws.Range("O26").Value = "9)  0 ,7 ,0 ,0 ,0 ,15, 0, 5, 0, 0"

# --- Row 17: new item "inanimate bones" + stats ---
$ws.Range("B17").Value = "inanimate bones"
$ws.Range("E17").Value = 3
$ws.Range("H17").Value = 11

# --- Row 18: new item "Shoulder pad" + stats ---
$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "Shoulder pad"
$ws.Range("E18").Value = 1
$ws.Range("F18").Value = 4
$ws.Range("H18").Value = 12
$ws.Range("M18").Value = 4

# --- Move the active selection (the "loading bar") to M18 ---
$ws.Range("M18").Select()
